# "Generate Report for Handback" - update the localization-status workbook
# after a successful handback for the de-de / zh-cn targets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdFile    = "92275e07-faf6-479a-a38a-950c2959146a.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16e0c8d940025351255dde198b7dc0af6c34f399/e2e/92275e07-faf6-479a-a38a-950c2959146a.md"
$zhHandbackXlf = "92275e07-faf6-479a-a38a-950c2959146a.4390e6f653466e5aead1aea3810d6008917612d5.zh-cn.xlf"
$deHandbackXlf = "92275e07-faf6-479a-a38a-950c2959146a.4390e6f653466e5aead1aea3810d6008917612d5.de-de.xlf"
$zhHandbackDate = "2016-09-04 07:05:18"
$deHandbackDate = "2016-09-04 07:05:26"

# ---- Overview sheet: status for both locales now shows "handed back" ----
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn sheet ----
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("I2").Value = $mdFile
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFile) | Out-Null
$zhcn.Range("J2").Value = $zhHandbackXlf
$zhcn.Range("K2").Value = $zhHandbackDate
$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet ----
$dede.Range("C2").Value = $newStatus
$dede.Range("I2").Value = $mdFile
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFile) | Out-Null
$dede.Range("J2").Value = $deHandbackXlf
$dede.Range("K2").Value = $deHandbackDate
$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
